# Update figure panel captions from "A:/B:/C:/D:" style to "(a)/(b)/(c)/(d)"
# style, and adjust each caption textbox's width to match the new
# (auto-fit) rendered size of the updated text (wrap="none" + spAutoFit).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# TextBox 104 (cNvPr id=105): "A: Genes" -> "(a) Genes"
$shA = $s.Shapes.Item(20)
$shA.TextFrame.TextRange.Text = "(a) Genes"
$shA.Width = 1098378 / 12700

# TextBox 105 (cNvPr id=106): "B: Aphid abundance" -> "(b) Aphid abundance"
$shB = $s.Shapes.Item(21)
$shB.TextFrame.TextRange.Text = "(b) Aphid abundance"
$shB.Width = 2076081 / 12700

# TextBox 108 (cNvPr id=109): "C: Bolting" -> "(c) Bolting"
$shC = $s.Shapes.Item(24)
$shC.TextFrame.TextRange.Text = "(c) Bolting"
$shC.Width = 1107996 / 12700

# TextBox 109 (cNvPr id=110): "D: Residuals" -> "(d) Residuals"
$shD = $s.Shapes.Item(25)
$shD.TextFrame.TextRange.Text = "(d) Residuals"
$shD.Width = 1391728 / 12700
